$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1 ("Version & History"): add a new change-history row (row 16)
# ---------------------------------------------------------------------------

# Clone formatting from the row above (row 15) so the new row gets the same
# styles (font/border/alignment) instead of the worksheet default.
$ws1.Range("A15:E15").Copy()
$ws1.Range("A16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Cells.Item(16, 1).Value = "V2.1"
$ws1.Cells.Item(16, 2).Value = "Fixed Camera signals by adding own lat-lon ego for the Camera`n"
$ws1.Cells.Item(16, 3).Value = "Bertalan Ádám"
$ws1.Cells.Item(16, 4).Value = 42825
$ws1.Cells.Item(16, 5).Value = "Draft version"

$ws1.Rows.Item(16).RowHeight = 42.75

# ---------------------------------------------------------------------------
# Sheet2 ("CommunicationMatrix"): fix Camera signals
# ---------------------------------------------------------------------------

# Row 20 ("Traffic sign meaning") was missing a Type entry - add it.
$ws2.Cells.Item(20, 10).Value = "PowertrainSystem"

# Row 21: own "Lateral EGO" signal for the Camera (copy of row 18's layout).
$ws2.Range("A18:J18").Copy()
$ws2.Range("A21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws2.Cells.Item(21, 1).Value = 17
$ws2.Cells.Item(21, 2).Value = "Lateral EGO"
$ws2.Cells.Item(21, 3).Value = 0
$ws2.Cells.Item(21, 4).Value = "-"
$ws2.Cells.Item(21, 5).Value = 1
$ws2.Cells.Item(21, 6).Value = "m"
$ws2.Cells.Item(21, 7).Value = "-"
$ws2.Cells.Item(21, 8).Value = "Output"
$ws2.Cells.Item(21, 9).Value = "Camera"
$ws2.Cells.Item(21, 10).Value = "PowertrainSystem"

# Row 22: own "Longitudinal EGO" signal for the Camera (copy of row 16's layout).
$ws2.Range("A16:J16").Copy()
$ws2.Range("A22:J22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws2.Cells.Item(22, 1).Value = 18
$ws2.Cells.Item(22, 2).Value = "Longitudinal EGO"
$ws2.Cells.Item(22, 3).Value = 0
$ws2.Cells.Item(22, 4).Value = "-"
$ws2.Cells.Item(22, 5).Value = 1
$ws2.Cells.Item(22, 6).Value = "m"
$ws2.Cells.Item(22, 7).Value = "-"
$ws2.Cells.Item(22, 8).Value = "Output"
$ws2.Cells.Item(22, 9).Value = "Camera"
$ws2.Cells.Item(22, 10).Value = "PowertrainSystem"

# ---------------------------------------------------------------------------
# View state: selections / scroll position / active sheet
# ---------------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("J5").Select()

$ws1.Activate()
$ws1.Range("D17").Select()
